# Reorders the comma-separated "Recorded By" names in column G of the
# "Session Analysis Results" sheet according to a fixed priority order:
#   backup@backdoor.com < System < dnasr281@gmail.com < admin@admin.com < system
# This mirrors the upstream change where the recorder list was re-sorted.

# NOTE: this PowerShell runtime compares strings (-eq/-ceq/-clike/switch, and
# hashtable key lookups) case-INsensitively even with "case sensitive"
# operators, which matters here because "System" and "system" are distinct,
# meaningful tokens. So do exact, case-sensitive comparisons by comparing
# character codes one by one instead of relying on string equality operators.
function Test-ExactMatch($a, $b) {
    if ($a.Length -ne $b.Length) { return $false }
    $ca = $a.ToCharArray()
    $cb = $b.ToCharArray()
    for ($i = 0; $i -lt $ca.Length; $i++) {
        if ([int]$ca[$i] -ne [int]$cb[$i]) { return $false }
    }
    return $true
}

function Get-RecorderRank($name) {
    if (Test-ExactMatch $name "backup@backdoor.com") { return 0 }
    if (Test-ExactMatch $name "System") { return 1 }
    if (Test-ExactMatch $name "dnasr281@gmail.com") { return 2 }
    if (Test-ExactMatch $name "admin@admin.com") { return 3 }
    if (Test-ExactMatch $name "system") { return 4 }
    return 999
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*"
    if ($parts.Count -le 1) {
        continue
    }

    # Build real objects carrying a precomputed rank, then sort by that
    # real property (sorting directly on an inline scriptblock's return
    # value was unreliable in this runtime).
    $items = @()
    foreach ($p in $parts) {
        $token = $p.Trim()
        $items += [PSCustomObject]@{ Token = $token; Rank = (Get-RecorderRank $token) }
    }

    $sortedItems = $items | Sort-Object -Property Rank

    $sortedTokens = @()
    foreach ($it in $sortedItems) {
        $sortedTokens += $it.Token
    }

    $newValue = [string]::Join(", ", $sortedTokens)

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
